$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44249
$ws.Range("K2").Value = 12000
$ws.Range("L2").Value = 12000
$ws.Range("M2").Value = 12000
$ws.Range("P2").Value = 667

$ws.Range("D3").Value = 44249
$ws.Range("I3").Value = 'Segunda'
$ws.Range("J3").Value = 200
$ws.Range("K3").Value = 10000
$ws.Range("L3").Value = 10000
$ws.Range("M3").Value = 10000
$ws.Range("P3").Value = 556

$ws.Range("D4").Value = 44726
$ws.Range("I4").Value = 'Primera'
$ws.Range("J4").Value = 300
$ws.Range("K4").Value = 14000
$ws.Range("L4").Value = 14000
$ws.Range("M4").Value = 14000
$ws.Range("P4").Value = 778

$ws.Range("D5").Value = 44614
$ws.Range("I5").Value = 'Primera'
$ws.Range("J5").Value = 300
$ws.Range("K5").Value = 15000
$ws.Range("L5").Value = 15000
$ws.Range("M5").Value = 15000
$ws.Range("N5").Value = '$/caja 18 kilos granel'
$ws.Range("P5").Value = 833

$ws.Range("D6").Value = 44229
$ws.Range("I6").Value = 'Primera'
$ws.Range("J6").Value = 200
$ws.Range("K6").Value = 15000
$ws.Range("L6").Value = 15000
$ws.Range("M6").Value = 15000
$ws.Range("P6").Value = 833

$ws.Range("D7").Value = 44719
$ws.Range("K7").Value = 15000
$ws.Range("L7").Value = 15000
$ws.Range("M7").Value = 15000
$ws.Range("P7").Value = 833

$ws.Range("D8").Value = 44628
$ws.Range("I8").Value = 'Primera'
$ws.Range("J8").Value = 300
$ws.Range("K8").Value = 15000
$ws.Range("L8").Value = 15000
$ws.Range("M8").Value = 15000
$ws.Range("P8").Value = 833

$ws.Range("D9").Value = 44596
$ws.Range("J9").Value = 150
$ws.Range("K9").Value = 14000
$ws.Range("L9").Value = 14000
$ws.Range("M9").Value = 14000
$ws.Range("P9").Value = 778

$ws.Range("D10").Value = 44238
$ws.Range("I10").Value = 'Primera'
$ws.Range("J10").Value = 300

$ws.Range("D11").Value = 44238
$ws.Range("I11").Value = 'Segunda'
$ws.Range("K11").Value = 10000
$ws.Range("L11").Value = 10000
$ws.Range("M11").Value = 10000
$ws.Range("P11").Value = 556

$ws.Range("D12").Value = 44238
$ws.Range("I12").Value = 'Tercera'
$ws.Range("J12").Value = 50
$ws.Range("K12").Value = 8000
$ws.Range("L12").Value = 8000
$ws.Range("M12").Value = 8000
$ws.Range("O12").Value = 'Provincia de Limarí'
$ws.Range("P12").Value = 444

$ws.Range("D13").Value = 44585
$ws.Range("J13").Value = 200
$ws.Range("K13").Value = 12000
$ws.Range("L13").Value = 12000
$ws.Range("M13").Value = 12000
$ws.Range("P13").Value = 667

$ws.Range("D14").Value = 44391
$ws.Range("J14").Value = 400

$ws.Range("D15").Value = 44721
$ws.Range("J15").Value = 300
$ws.Range("K15").Value = 15000
$ws.Range("L15").Value = 15000
$ws.Range("M15").Value = 15000
$ws.Range("P15").Value = 833

$ws.Range("D16").Value = 44627
$ws.Range("I16").Value = 'Primera'
$ws.Range("J16").Value = 300
$ws.Range("K16").Value = 15000
$ws.Range("L16").Value = 15000
$ws.Range("M16").Value = 15000
$ws.Range("P16").Value = 833

$ws.Range("D17").Value = 44635
$ws.Range("I17").Value = 'Primera'
$ws.Range("J17").Value = 300
$ws.Range("K17").Value = 15000
$ws.Range("L17").Value = 15000
$ws.Range("M17").Value = 15000
$ws.Range("P17").Value = 833

$ws.Range("D18").Value = 44243
$ws.Range("I18").Value = 'Especial'
$ws.Range("K18").Value = 12000
$ws.Range("L18").Value = 12000
$ws.Range("M18").Value = 12000
$ws.Range("P18").Value = 667

$ws.Range("D19").Value = 44243
$ws.Range("K19").Value = 10000
$ws.Range("L19").Value = 10000
$ws.Range("M19").Value = 10000
$ws.Range("P19").Value = 556

$ws.Range("D20").Value = 44243
$ws.Range("I20").Value = 'Segunda'
$ws.Range("J20").Value = 150
$ws.Range("K20").Value = 8000
$ws.Range("L20").Value = 8000
$ws.Range("M20").Value = 8000
$ws.Range("N20").Value = '$/bandeja 18 kilos'
$ws.Range("P20").Value = 444

$ws.Range("D21").Value = 44631
$ws.Range("J21").Value = 300
$ws.Range("K21").Value = 15000
$ws.Range("L21").Value = 15000
$ws.Range("M21").Value = 15000
$ws.Range("P21").Value = 833

$ws.Range("D22").Value = 44235
$ws.Range("J22").Value = 400
$ws.Range("K22").Value = 13000
$ws.Range("L22").Value = 13000
$ws.Range("M22").Value = 13000
$ws.Range("P22").Value = 722

$ws.Range("D23").Value = 44235
$ws.Range("I23").Value = 'Segunda'
$ws.Range("J23").Value = 200
$ws.Range("K23").Value = 11000
$ws.Range("L23").Value = 11000
$ws.Range("M23").Value = 11000
$ws.Range("P23").Value = 611

$ws.Range("D24").Value = 44235
$ws.Range("I24").Value = 'Tercera'
$ws.Range("J24").Value = 100
$ws.Range("K24").Value = 9000
$ws.Range("L24").Value = 9000
$ws.Range("M24").Value = 9000
$ws.Range("P24").Value = 500

$ws.Range("D25").Value = 44725
$ws.Range("J25").Value = 400
$ws.Range("K25").Value = 15000
$ws.Range("L25").Value = 15000
$ws.Range("M25").Value = 15000
$ws.Range("O25").Value = 'Región de Arica y Parinacota'
$ws.Range("P25").Value = 833

$ws.Range("D26").Value = 44630
$ws.Range("J26").Value = 300
$ws.Range("K26").Value = 15000
$ws.Range("L26").Value = 15000
$ws.Range("M26").Value = 15000
$ws.Range("P26").Value = 833

$ws.Range("D27").Value = 44383
$ws.Range("K27").Value = 16000
$ws.Range("L27").Value = 16000
$ws.Range("M27").Value = 16000
$ws.Range("P27").Value = 889

$ws.Range("D28").Value = 44383
$ws.Range("I28").Value = 'Segunda'
$ws.Range("J28").Value = 200
$ws.Range("K28").Value = 12000
$ws.Range("L28").Value = 12000
$ws.Range("M28").Value = 12000
$ws.Range("P28").Value = 667

$ws.Range("D29").Value = 44396
$ws.Range("J29").Value = 250
$ws.Range("K29").Value = 15000
$ws.Range("L29").Value = 15000
$ws.Range("M29").Value = 15000
$ws.Range("P29").Value = 833

$ws.Range("D30").Value = 44396
$ws.Range("J30").Value = 150
$ws.Range("K30").Value = 12000
$ws.Range("L30").Value = 12000
$ws.Range("M30").Value = 12000
$ws.Range("P30").Value = 667

$ws.Range("D31").Value = 44245
$ws.Range("J31").Value = 300

$ws.Range("D32").Value = 44245
